# error solve ifrs list
# Corrects the financial figures for the actual-year columns (2014-2018,
# rows 2-6) and clears out the erroneous estimate-year columns
# (2019E-2021E, rows 7-9) that should no longer carry data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (2014/12 IFRS연결) ---
$ws.Range("D2").Value = 1359
$ws.Range("E2").Value = 72
$ws.Range("F2").Value = 72
$ws.Range("G2").Value = 153
$ws.Range("H2").Value = 159
$ws.Range("I2").Value = 136
$ws.Range("J2").ClearContents()
$ws.Range("K2").Value = 4056
$ws.Range("L2").Value = 1300
$ws.Range("M2").Value = 2756
$ws.Range("N2").Value = 2603
$ws.Range("O2").ClearContents()
$ws.Range("P2").Value = 48
$ws.Range("Q2").Value = 61
$ws.Range("R2").Value = -406
$ws.Range("S2").Value = 294
$ws.Range("T2").Value = 330
$ws.Range("U2").Value = -269
$ws.Range("V2").Value = 879
$ws.Range("W2").Value = 5.29
$ws.Range("X2").Value = 11.69
$ws.Range("Y2").Value = 7.63
$ws.Range("Z2").Value = 4.45
$ws.Range("AA2").Value = 47.18
$ws.Range("AB2").Value = 6214.27
$ws.Range("AC2").Value = 1690
$ws.Range("AD2").Value = 26.86
$ws.Range("AE2").Value = 27237
$ws.Range("AF2").Value = 1.67
$ws.Range("AG2").Value = 240
$ws.Range("AH2").Value = 0.53
$ws.Range("AI2").Value = 16.91
$ws.Range("AJ2").Value = 9603921

# --- Row 3 (2015/12 IFRS연결) ---
$ws.Range("D3").Value = 1999
$ws.Range("E3").Value = 113
$ws.Range("F3").Value = 113
$ws.Range("G3").Value = 119
$ws.Range("H3").Value = 103
$ws.Range("I3").Value = 80
$ws.Range("J3").ClearContents()
$ws.Range("K3").Value = 5446
$ws.Range("L3").Value = 2192
$ws.Range("M3").Value = 3254
$ws.Range("N3").Value = 2793
$ws.Range("O3").ClearContents()
$ws.Range("P3").Value = 48
$ws.Range("Q3").Value = 239
$ws.Range("R3").Value = -1407
$ws.Range("S3").Value = 1190
$ws.Range("T3").Value = 1022
$ws.Range("U3").Value = -783
$ws.Range("V3").Value = 1545
$ws.Range("W3").Value = 5.63
$ws.Range("X3").Value = 5.16
$ws.Range("Y3").Value = 2.97
$ws.Range("Z3").Value = 2.17
$ws.Range("AA3").Value = 67.34999999999999
$ws.Range("AB3").Value = 6600.51
$ws.Range("AC3").Value = 833
$ws.Range("AD3").Value = 82.45
$ws.Range("AE3").Value = 29219
$ws.Range("AF3").Value = 2.35
$ws.Range("AG3").Value = 130
$ws.Range("AH3").Value = 0.19
$ws.Range("AI3").Value = 15.53
$ws.Range("AJ3").Value = 9603921

# --- Row 4 (2016/12 IFRS연결) ---
$ws.Range("D4").Value = 2669
$ws.Range("E4").Value = 202
$ws.Range("F4").Value = 202
$ws.Range("G4").Value = 285
$ws.Range("H4").Value = 245
$ws.Range("I4").Value = 154
$ws.Range("J4").Value = 91
$ws.Range("K4").Value = 6707
$ws.Range("L4").Value = 3241
$ws.Range("M4").Value = 3466
$ws.Range("N4").Value = 2902
$ws.Range("O4").Value = 564
$ws.Range("P4").Value = 48
$ws.Range("Q4").Value = 128
$ws.Range("R4").Value = -656
$ws.Range("S4").Value = 739
$ws.Range("T4").Value = 805
$ws.Range("U4").Value = -677
$ws.Range("V4").Value = 2295
$ws.Range("W4").Value = 7.55
$ws.Range("X4").Value = 9.18
$ws.Range("Y4").Value = 5.39
$ws.Range("Z4").Value = 4.03
$ws.Range("AA4").Value = 93.48
$ws.Range("AB4").Value = 6886.41
$ws.Range("AC4").Value = 1599
$ws.Range("AD4").Value = 23.39
$ws.Range("AE4").Value = 30363
$ws.Range("AF4").Value = 1.23
$ws.Range("AG4").Value = 200
$ws.Range("AH4").Value = 0.53
$ws.Range("AI4").Value = 12.45
$ws.Range("AJ4").Value = 9603921

# --- Row 5 (2017/12 IFRS연결) ---
$ws.Range("D5").Value = 2983
$ws.Range("E5").Value = 101
$ws.Range("F5").Value = 101
$ws.Range("G5").Value = 86
$ws.Range("H5").Value = 58
$ws.Range("I5").Value = 38
$ws.Range("J5").Value = 20
$ws.Range("K5").Value = 7408
$ws.Range("L5").Value = 3971
$ws.Range("M5").Value = 3437
$ws.Range("N5").Value = 2888
$ws.Range("O5").Value = 548
$ws.Range("P5").Value = 48
$ws.Range("Q5").Value = 109
$ws.Range("R5").Value = -814
$ws.Range("S5").Value = 586
$ws.Range("T5").Value = 743
$ws.Range("U5").Value = -634
$ws.Range("V5").Value = 2890
$ws.Range("W5").Value = 3.4
$ws.Range("X5").Value = 1.96
$ws.Range("Y5").Value = 1.33
$ws.Range("Z5").Value = 0.83
$ws.Range("AA5").Value = 115.56
$ws.Range("AB5").Value = 6885.68
$ws.Range("AC5").Value = 400
$ws.Range("AD5").Value = 82.19
$ws.Range("AE5").Value = 30218
$ws.Range("AF5").Value = 1.09
$ws.Range("AG5").Value = 100
$ws.Range("AH5").Value = 0.3
$ws.Range("AI5").Value = 24.86
$ws.Range("AJ5").Value = 9603921

# --- Row 6 (2018/12 IFRS연결) ---
$ws.Range("D6").Value = 4383
$ws.Range("E6").Value = 241
$ws.Range("F6").Value = 241
$ws.Range("G6").Value = 189
$ws.Range("H6").Value = 126
$ws.Range("I6").Value = 95
$ws.Range("K6").Value = 8798
$ws.Range("L6").Value = 5200
$ws.Range("M6").Value = 3598
$ws.Range("N6").Value = 2990
$ws.Range("P6").Value = 48
$ws.Range("Q6").Value = -14
$ws.Range("R6").Value = -730
$ws.Range("S6").Value = 771
$ws.Range("T6").Value = 525
$ws.Range("U6").Value = -539
$ws.Range("V6").Value = 3815
$ws.Range("W6").Value = 5.51
$ws.Range("X6").Value = 2.87
$ws.Range("Y6").Value = 3.23
$ws.Range("Z6").Value = 1.55
$ws.Range("AA6").Value = 144.5
$ws.Range("AB6").Value = 7042.74
$ws.Range("AC6").Value = 990
$ws.Range("AD6").Value = 19.76
$ws.Range("AE6").Value = 31281
$ws.Range("AF6").Value = 0.62
$ws.Range("AG6").Value = 140
$ws.Range("AH6").Value = 0.72
$ws.Range("AI6").Value = 14.08
$ws.Range("AJ6").Value = 9603921

# --- Rows 7-9 (2019/12(E), 2020/12(E), 2021/12(E)) ---
# These estimate rows no longer carry any financial figures - only the
# leading A/B/C identifier columns remain.
$ws.Range("D7:AJ9").ClearContents()
